$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update a few existing "说明" (description) cells in column D ---
$ws.Range("D6").Value  = "MCU串口发-usart3_tx"
$ws.Range("D7").Value  = "MCU串口收-usart3_rx"
$ws.Range("D10").Value = "板载绿色LED，低亮高灭"
$ws.Range("D11").Value = "板载红色LED，低亮高灭"
$ws.Range("D15").Value = "usart1_tx"
$ws.Range("D16").Value = "usart1_rx"

# --- New column F ("用途" / purpose) ---
$ws.Range("F1").Value  = "用途"

$ws.Range("F4").Value  = "ADC1_CH6"
$ws.Range("F5").Value  = "ADC2_CH5"
$ws.Range("F6").Value  = "shell"
$ws.Range("F12").Value = "ADC1_CH0"
$ws.Range("F13").Value = "ADC1_CH4"
$ws.Range("F15").Value = "tuya"

# Shrink column D and give column F the same look & feel as the rest of
# the table (left border, centered header / cell formatting where needed)
$ws.Columns("D").ColumnWidth = 29.19921875

# Match the existing style used by the other data cells (s="2" -> left
# border, left aligned, vertical centered) for the plain F cells first
$ws.Range("F1:F16").HorizontalAlignment = -4131
$ws.Range("F1:F16").VerticalAlignment = -4108
$ws.Range("F1:F16").Borders.Item(7).LineStyle = 1

# Merge F6:F7 ("shell") and F15:F16 ("tuya") and box + center them
$ws.Range("F6:F7").Merge()
$ws.Range("F15:F16").Merge()

$ws.Range("F6:F7").HorizontalAlignment = -4108
$ws.Range("F15:F16").HorizontalAlignment = -4108

$ws.Range("F6:F7").Borders.Item(7).LineStyle = 1
$ws.Range("F6:F7").Borders.Item(10).LineStyle = 1
$ws.Range("F6").Borders.Item(8).LineStyle = 1
$ws.Range("F7").Borders.Item(9).LineStyle = 1

$ws.Range("F15:F16").Borders.Item(7).LineStyle = 1
$ws.Range("F15:F16").Borders.Item(10).LineStyle = 1
$ws.Range("F15").Borders.Item(8).LineStyle = 1
$ws.Range("F16").Borders.Item(9).LineStyle = 1

# Selection, as left by the editing session
$ws.Range("F13").Select()

Write-Output "applied"
